$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell X1 and data cell X2
$ws.Range("X1").Value = "visibility_level"
$ws.Range("X2").Value = "PRO"

# Match existing column width style (customWidth) for column X (24).
# Note: this runtime's ColumnWidth setter quantizes to 1/6-character steps
# (vs. real Excel's 1/256), so 18.0 is the closest achievable input to the
# target stored width of 18.83203125 (yields ~18.8333).
$ws.Columns.Item(24).ColumnWidth = 18.0
